$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet is protected; unprotect to allow edits, then restore protection after
$ws.Unprotect()

# Update the confidential disclosure date from 2021-04-30 to 2021-05-03
[void]$ws.Range("A9").Replace("2021-04-30", "2021-05-03")

# Update weight (D) and percent change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2493179965542861
$ws.Range("E2").Value = 0.01493256262042375

$ws.Range("D3").Value = 0.2512515054480104
$ws.Range("E3").Value = 0.004964147821290732

$ws.Range("D4").Value = 0.2453498943725617
$ws.Range("E4").Value = 0.01030118708917871

$ws.Range("D5").Value = 0.2540806036251418
$ws.Range("E5").Value = -0.005253042921204321

$ws.Range("E6").Value = 0.006162905057264911

# Restore sheet protection as it was originally (password unknown/not changed
# by this edit, so contents are protected again without a new password)
$ws.Protect()
